$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new blank column before column A (shifts old A/B/C -> B/C/D)
# ---------------------------------------------------------------------------
$ws.Columns("A").Insert()

# Width of the new column A (closest achievable value given the runtime's
# character-width snapping; target stored width is 4.42578125).
$ws.Columns("A").ColumnWidth = 3.71

# ---------------------------------------------------------------------------
# 2. Highlight certain rows in the new column A with a yellow fill
#    (these cells stay empty, just colored)
# ---------------------------------------------------------------------------
$yellow = 65535
$highlightRows = @(2, 4, 9, 41, 42, 45, 46, 47)
foreach ($r in $highlightRows) {
    $ws.Cells.Item($r, 1).Interior.Color = $yellow
}

# ---------------------------------------------------------------------------
# 3. Add the new "Subtitle" text (column D) and mark the related "Code"
#    cells (column B) with a green fill for the newly recorded / highlighted
#    voice lines.
# ---------------------------------------------------------------------------
$green = 5287936

$ws.Cells.Item(41, 4).Value = "This is my home… they can’t just take this place away from me"
$ws.Cells.Item(41, 2).Interior.Color = $green

$ws.Cells.Item(42, 4).Value = "He should understand… he grew up here with us. I can’t just leave, not like he did."
$ws.Cells.Item(42, 2).Interior.Color = $green

$ws.Cells.Item(45, 4).Value = "Alex used to be stuck in these books of hers for hours, she was always the smart one."
$ws.Cells.Item(45, 2).Interior.Color = $green

$ws.Cells.Item(46, 4).Value = "That tool head is loose on that hoe, I should wedge some wood in there later…"
$ws.Cells.Item(46, 2).Interior.Color = $green

$ws.Cells.Item(47, 4).Value = "Thar she is… the beauty. I remember taking Nichola out on that old boat. Taught him to fish like my paps taught me."
$ws.Cells.Item(47, 2).Interior.Color = $green
